$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fn1", "Tnfrsf11b", "ECs", 2, 1, 40.7349555, 81.469911, 0.05567871843833241, 0.03826666865920979, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 4.192686029793, 16.770744119172, 0.002127821065063995, 0.001042452480425824),
    @("ECs", "Fn1", "Tnfrsf11b", "FAPs", 2, 1, 40.7349555, 81.469911, 0.05567871843833241, 0.03826666865920979, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 88.392708060677, 530.356248364062, 0.04485999306245371, 0.03296640761362634),
    @("ECs", "Fn1", "Tnfrsf11b", "sCs", 2, 1, 40.7349555, 81.469911, 0.05567871843833241, 0.03826666865920979, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 17.12466978003375, 68.498679120135, 0.008690904310814717, 0.00425780856515762),
    @("FAPs", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 350.3919066666667, 1051.17572, 0.4789344206933965, 0.4937404802104949, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 36.06443738557333, 216.38662431344, 0.01830298501328351, 0.01345037355788202),
    @("FAPs", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 350.3919066666667, 1051.17572, 0.4789344206933965, 0.4937404802104949, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 760.3319834929156, 6842.987851436241, 0.3858744488429991, 0.4253531989137333),
    @("FAPs", "Fn1", "Tnfrsf11b", "sCs", 3, 1, 350.3919066666667, 1051.17572, 0.4789344206933965, 0.4937404802104949, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 147.3021296233667, 883.8127777401999, 0.07475698683711392, 0.05493690773887956),
    @("M1", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 243.8287033333334, 731.48611, 0.3332781281688242, 0.3435812836494235, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 25.09631311928667, 150.57787871572, 0.01273657586835725, 0.009359768537939576),
    @("M1", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 243.8287033333334, 731.48611, 0.3332781281688242, 0.3435812836494235, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 529.0954445882912, 4761.859001294621, 0.2685200905635067, 0.2959923359430933),
    @("M1", "Fn1", "Tnfrsf11b", "sCs", 3, 1, 243.8287033333334, 731.48611, 0.3332781281688242, 0.3435812836494235, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 102.5037581660583, 615.02254899635, 0.05202146173696028, 0.03822917916839052),
    @("M2", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 48.737294, 146.211882, 0.06661674320651284, 0.06867617226847689, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 5.016334722244, 30.098008333464, 0.002545829240623991, 0.00187086173518255),
    @("M2", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 48.737294, 146.211882, 0.06661674320651284, 0.06867617226847689, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 105.7573611491827, 951.8162503426441, 0.05367269078574952, 0.05916393476810643),
    @("M2", "Fn1", "Tnfrsf11b", "sCs", 3, 1, 48.737294, 146.211882, 0.06661674320651284, 0.06867617226847689, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 20.488792867895, 122.93275720737, 0.01039822318013933, 0.007641375765187903),
    @("Neutro", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 22.832077, 68.496231, 0.03120810544755168, 0.03217289111905, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 2.350014357302, 14.100086143812, 0.001192650729660504, 0.0008764470837064035),
    @("Neutro", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 22.832077, 68.496231, 0.03120810544755168, 0.03217289111905, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 49.54440460061133, 445.899641405502, 0.02514417416809032, 0.02771667040538572),
    @("Neutro", "Fn1", "Tnfrsf11b", "sCs", 3, 1, 22.832077, 68.496231, 0.03120810544755168, 0.03217289111905, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 9.598433930222498, 57.59060358133499, 0.004871280549800859, 0.003579773629957874),
    @("sCs", "Fn1", "Tnfrsf11b", "ECs", 2, 1, 25.082339, 50.164678, 0.03428388404538221, 0.02356250409334498, 1, 0.5, 0.102926, 0.205852, 0.03821605677617539, 0.02724178813968781, 2.581624823914, 10.326499295656, 0.00131019485918614, 0.0006418847445514305),
    @("sCs", "Fn1", "Tnfrsf11b", "FAPs", 2, 1, 25.082339, 50.164678, 0.03428388404538221, 0.02356250409334498, 3, 1, 2.169947333333333, 6.509842000000001, 0.8056937070514454, 0.8614914433031574, 54.42735462681267, 326.564127760876, 0.0276223096286459, 0.02029889565921232),
    @("sCs", "Fn1", "Tnfrsf11b", "sCs", 2, 1, 25.082339, 50.164678, 0.03428388404538221, 0.02356250409334498, 2, 1, 0.4203924999999999, 0.8407849999999999, 0.1560902361723793, 0.1112667685571547, 10.5444271980575, 42.17770879223, 0.005351379557550176, 0.002621723689581226)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
